$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 gets the "Link" header text (was an empty inline string cell before).
$ws.Cells.Item(1, 4).Value = "Link"

# New column D (rows 2-16): product page links.
$links = @{
    2  = "https://sneakerhead.ru/shoes/sneakers/slipstream-hi-xtreme-cordura-39327202/"
    3  = "https://sneakerhead.ru/shoes/sneakers/slipstream-hi-xtreme-cordura-39327201/"
    4  = "https://sneakerhead.ru/shoes/sneakers/slipstream-hi-xtreme-39327102/"
    5  = "https://sneakerhead.ru/shoes/sneakers/zig-kinetica-2-5-edge-winter-100073985/"
    6  = "https://sneakerhead.ru/shoes/sneakers/wmns-1-acclimate-DC7723-401/"
    7  = "https://sneakerhead.ru/shoes/sneakers/wmns-1-acclimate-DC7723-001/"
    8  = "https://sneakerhead.ru/shoes/sneakers/wave-mujin-tl-gtx-D1GA237302/"
    9  = "https://sneakerhead.ru/shoes/sneakers/wave-mujin-tl-gtx-D1GA237301/"
    10 = "https://sneakerhead.ru/shoes/boots/jasper-boots-HK-1323-010/"
    11 = "https://sneakerhead.ru/shoes/boots/jasper-boots-HK-1323-008/"
    12 = "https://sneakerhead.ru/shoes/boots/jasper-boots-HK-1323-012/"
    13 = "https://sneakerhead.ru/shoes/boots/jasper-boots-HK-1323-006/"
    14 = "https://sneakerhead.ru/shoes/sneakers/wmns-1-acclimate-DC7723-100/"
    15 = "https://sneakerhead.ru/shoes/sneakers/wmns-1-acclimate-DC7723-500/"
    16 = "https://sneakerhead.ru/shoes/sneakers/cl-lthr-mid-gtx-thin-GZ6887/"
}

foreach ($row in $links.Keys) {
    $ws.Cells.Item($row, 4).Value = $links[$row]
}

# E1 becomes a new (empty) inline-string cell, extending the used range
# from A1:D16 to A1:E16, without actually giving it a value/format.
$ws.Cells.Item(1, 5).Font.Bold = $false
